$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.352.83"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "1.592.03"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").Value = "1.631.05"
$ws.Range("E13").Value = "  +3.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.06"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "26.353.42"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "212.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.48%  "

$ws.Range("E21").Value = "  -0.31%  "

$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("E23").Value = "  -1.23%  "

$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("E26").Value = "  -0.23%  "

$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("E31").Value = "  +1.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "

$ws.Range("E33").Value = "  +1.21%  "

$ws.Range("D34").Value = "1.343.27"
$ws.Range("E34").Value = "  +4.40%  "

$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.603"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -16.72%  "

$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.30%  "

$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("E44").Value = "  -0.60%  "

$ws.Range("D45").Value = "1.728.29"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +3.37%  "

$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0987"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.24%  "

$ws.Range("E51").Value = "  -0.86%  "
